# Update column F (dSF) values after repulling data / recalculating mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    4  = -1
    5  = 1
    6  = -1
    7  = -4
    8  = 3
    9  = 5
    10 = 2
    11 = -3
    12 = 3
    13 = 1
    14 = -1
    15 = 6
    17 = 4
    18 = -3
    19 = 2
    20 = -4
    22 = -1
    25 = -1
    26 = 5
    27 = -3
    29 = -2
    30 = -6
    31 = 1
    32 = -2
    33 = -2
    34 = 5
    35 = 3
    36 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
